$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Toggle Runmode from Y to N for all suites except "C Suite" (Authoring module, row 4)
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Update the selected/active cell shown in the sheet view
$ws.Range("C4").Select()
